$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "price delta per node" ratios on the existing table (col G, rows 4,6,8,10,12) ---
$ws.Range("G4").Formula  = "=36.3/33.6"
$ws.Range("G6").Formula  = "=29.9/22.1"
$ws.Range("G8").Formula  = "=36.5/13.4"
$ws.Range("G10").Formula = "=32.4/9.81"
$ws.Range("G12").Formula = "=20.4/9.94"

# --- New column widths ---
$ws.Columns.Item(1).ColumnWidth = 19.498697916666668
$ws.Columns.Item(5).ColumnWidth = 16.830729166666668
$ws.Columns.Item(6).ColumnWidth = 18.830729166666668
$ws.Columns.Item(7).ColumnWidth = 20.166666666666668

# --- New "performance / price" analysis block (rows 22-30) ---
# Shared-string table order matches the author's original entry order:
# #Compute node, Perf Delta, Total partner price diff, Total non-partner price diff,
# Simulation time, Price delta / node (Partner), Price delta / node (NonPartner).
$ws.Range("A23").Value = "#Compute node"
$ws.Range("D23").Value = "Perf Delta"
$ws.Range("F23").Value = "Total partner price diff"
$ws.Range("H23").Value = "Total non-partner price diff"

$ws.Range("B22:C22").Merge() | Out-Null
$ws.Range("B22").Value = "Simulation time"
$ws.Range("B22:C22").HorizontalAlignment = -4108

$ws.Range("E23").Value = "Price delta / node (Partner)"
$ws.Range("G23").Value = "Price delta / node (NonPartner)"

$ws.Range("B23").Value = "AWS-EC2"
$ws.Range("C23").Value = "ARCHER"

$ws.Range("A24").Value = 1
$ws.Range("B24").Value = 36.3
$ws.Range("C24").Value = 33.6
$ws.Range("D24").Formula = "=B24/C24"
$ws.Range("E24").Value = 7.55
$ws.Range("F24").Formula = "=E24*D24"
$ws.Range("G24").Value = 3.14
$ws.Range("H24").Formula = "=D24*G24"

$ws.Range("A25").Value = 2
$ws.Range("B25").Value = 29.9
$ws.Range("C25").Value = 22.1
$ws.Range("D25").Formula = "=B25/C25"
$ws.Range("E25").Value = 7.55
$ws.Range("F25").Formula = "=E25*D25"
$ws.Range("G25").Value = 3.14
$ws.Range("H25").Formula = "=D25*G25"

$ws.Range("A26").Value = 3
$ws.Range("B26").Value = 36.5
$ws.Range("C26").Value = 13.4
$ws.Range("D26").Formula = "=B26/C26"
$ws.Range("E26").Value = 7.55
$ws.Range("F26").Formula = "=E26*D26"
$ws.Range("G26").Value = 3.14
$ws.Range("H26").Formula = "=D26*G26"

$ws.Range("A27").Value = 4
$ws.Range("B27").Value = 32.4
$ws.Range("C27").Value = 9.81
$ws.Range("D27").Formula = "=B27/C27"
$ws.Range("E27").Value = 7.55
$ws.Range("F27").Formula = "=E27*D27"
$ws.Range("G27").Value = 3.14
$ws.Range("H27").Formula = "=D27*G27"

$ws.Range("A28").Value = 5
$ws.Range("B28").Value = 20.4
$ws.Range("C28").Value = 9.94
$ws.Range("D28").Formula = "=B28/C28"
$ws.Range("E28").Value = 7.55
$ws.Range("F28").Formula = "=E28*D28"
$ws.Range("G28").Value = 3.14
$ws.Range("H28").Formula = "=D28*G28"

$ws.Range("C29").Value = 9.64
$ws.Range("C30").Value = 18.6

# --- Move / resize the chart to make room for the new columns ---
# (computed after the column-width changes above, since those shift the
# pixel offsets of every column at/after the resized ones)
$co = $ws.ChartObjects().Item(1)
$co.Left = 766.1875
$co.Top = 26.0
$co.Width = 612.8125
$co.Height = 561.0

# --- Update selection to match the author's final cursor position ---
$ws.Range("G24").Select() | Out-Null
